# Weekly update: a new price-survey week is inserted right before the
# existing row 182 of the "Terminal La Palmera de La Serena - Espinaca"
# data block, pushing every subsequent record down by one row (old row
# 237 becomes row 238). Only the new row's Fecha/Volumen/Precio columns
# carry fresh data; every other column mirrors the constant values used
# throughout the rest of the block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row above row 182; this shifts rows 182-237 down to
# 183-238 (including their formatting), exactly like a manual Excel
# "Insert Row" above the first record of the new week.
$ws.Rows.Item(182).Insert()

# Populate the newly inserted row 182 with this week's record.
$ws.Range("A182").Value = 8
$ws.Range("B182").Value = "Terminal La Palmera de La Serena"
$ws.Range("C182").Value = "Coquimbo"
$ws.Range("D182").Value = 44627
$ws.Range("E182").Value = 4
$ws.Range("F182").Value = 100112012
$ws.Range("G182").Value = "Espinaca"
$ws.Range("H182").Value = "Sin especificar"
$ws.Range("I182").Value = "Primera"
$ws.Range("J182").Value = 2000
$ws.Range("K182").Value = 500
$ws.Range("L182").Value = 600
$ws.Range("M182").Value = 550
$ws.Range("N182").Value = "$/atado 300 a 500 gramos"
$ws.Range("O182").Value = "Provincia del Elquí"
$ws.Range("P182").Value = 1100
$ws.Range("Q182").Value = 0.5
$ws.Range("R182").Value = "Hortaliza"
